$wb = $excel.ActiveWorkbook

# "Test Results" sheet: mark the Create and Read tests as passed (TRUE)
# for every existing result row (rows 2-24), matching the commit
# "Added Create and Read tests to project".
$wsResults = $wb.Worksheets.Item("Test Results")
$wsResults.Range("B2:C24").Value = $true

# Reflect the author's last cell selection on the "Device" sheet, which
# moved from C38 to D36 while the workbook was being edited.
$wsDevice = $wb.Worksheets.Item("Device")
$wsDevice.Range("D36").Select()

# "Test Results" stays the active/visible sheet in the saved workbook,
# so re-activate it after touching the "Device" sheet's selection.
$wsResults.Activate()
